# Add a new "Phase 5" worksheet (censoring throughout pregnancy), inserted
# right after "Phase4" and before "OLD".

$wb = $excel.ActiveWorkbook

# Phase1's view now has the whole of column A selected (and is no longer
# scrolled down to row 16).
$phase1 = $wb.Worksheets.Item("Phase1")
$phase1.Activate() | Out-Null
$phase1.Columns.Item(1).Select() | Out-Null

$phase4 = $wb.Worksheets.Item("Phase4")
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $phase4)
$newSheet.Name = "Phase 5"

# Headers
$newSheet.Cells.Item(1, 1).Value = "gestweek_conception"
$newSheet.Cells.Item(1, 2).Value = "p_censoring"

# Data rows: gestweek_conception 0..40 in column A, p_censoring in column B.
# Weeks 0-5 -> 0, weeks 6-40 -> 0.02
for ($week = 0; $week -le 40; $week++) {
    $row = $week + 2
    $newSheet.Cells.Item($row, 1).Value = $week
    if ($week -le 5) {
        $newSheet.Cells.Item($row, 2).Value = 0
    } else {
        $newSheet.Cells.Item($row, 2).Value = 0.02
    }
}

# Column widths (raw stored widths end up ~0.8333 wider than the ColumnWidth
# value supplied, so back the numbers off to land on the target widths of
# 19 and 16.5).
$newSheet.Columns.Item(1).ColumnWidth = 18.16666666666667
$newSheet.Columns.Item(2).ColumnWidth = 15.666666666666666

# Make "Phase 5" the active sheet/tab, scrolled so row 13 is at the top with
# B8:B42 selected.
$newSheet.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$newSheet.Range("B8:B42").Select() | Out-Null

# "Phase4" is no longer the tab-selected sheet now that Phase 5 is active/selected.
